$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1595
$ws1.Range("F4").Value = 2116
$ws1.Range("F5").Value = 9121
$ws1.Range("F6").Value = 269
$ws1.Range("F8").Value = 1271
$ws1.Range("F10").Value = 632
$ws1.Range("F14").Value = 296
$ws1.Range("F17").Value = 1505
$ws1.Range("F19").Value = 570
$ws1.Range("F20").Value = 53
$ws1.Range("F21").Value = 1383
$ws1.Range("F22").Value = 83
$ws1.Range("F24").Value = 16
$ws1.Range("F26").Value = 65
$ws1.Range("F27").Value = 66
$ws1.Range("F28").Value = 314
$ws1.Range("F29").Value = 314
$ws1.Range("F30").Value = 1072
$ws1.Range("F33").Value = 230
$ws1.Range("F34").Value = 201
$ws1.Range("F35").Value = 61
$ws1.Range("F37").Value = 609
$ws1.Range("F40").Value = 75
$ws1.Range("F42").Value = 57
$ws1.Range("F43").Value = 497
$ws1.Range("F45").Value = 687
$ws1.Range("F46").Value = 219
$ws1.Range("F48").Value = 46

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 55
$ws2.Range("F21").Value = 76
$ws2.Range("F26").Value = 1034
$ws2.Range("F27").Value = 232
$ws2.Range("F30").Value = 230
$ws2.Range("F32").Value = 152
$ws2.Range("F35").Value = 115
$ws2.Range("F42").Value = 20

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 751
$ws3.Range("F7").Value = 2093
$ws3.Range("F8").Value = 3157

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1595
$ws4.Range("F4").Value = 751
$ws4.Range("F5").Value = 9121
$ws4.Range("F9").Value = 269
$ws4.Range("F10").Value = 2093
$ws4.Range("F11").Value = 3157
$ws4.Range("F15").Value = 1271
$ws4.Range("F17").Value = 632
$ws4.Range("F20").Value = 296
$ws4.Range("F22").Value = 1505
$ws4.Range("F23").Value = 570
$ws4.Range("F24").Value = 1383
$ws4.Range("F28").Value = 314
$ws4.Range("F29").Value = 314
$ws4.Range("F32").Value = 230
$ws4.Range("F34").Value = 232
$ws4.Range("F35").Value = 61
$ws4.Range("F37").Value = 609
$ws4.Range("F39").Value = 230
$ws4.Range("F41").Value = 152
$ws4.Range("F42").Value = 498
$ws4.Range("F43").Value = 687
$ws4.Range("F46").Value = 219
$ws4.Range("F50").Value = 20
